$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G7").Value = "Pass"
$ws.Range("G12").Value = "Pass"
$ws.Range("C2").Select()
